# Remove the "TEST POINT" test row (row 27) that was used for a temporary
# function/style test, shifting the rows below it up, then move the
# selection to where editing continued (B29), matching the new layout
# after the row was removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(27).Delete()

$ws.Range("B29").Select()

$wb.Save()
